# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 4902
$wsExhibit.Range("F5").Value = 24
$wsExhibit.Range("F6").Value = 25
$wsExhibit.Range("F8").Value = 489

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 4902
$wsAll.Range("F6").Value = 24
$wsAll.Range("F7").Value = 25
$wsAll.Range("F10").Value = 489
